# Update the COVID-19 countries data sheet ("Pais") with refreshed figures
# and re-sort the data body by total cases (column B) descending, mirroring
# the original data-refresh workflow used to build this workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- 1. Update the "last updated" timestamp in the title row ---------------
$ws.Range("A1").Value = "Datos actualizados a 13 de Mayo de 2020 a las 16:05"

# --- 2. Refresh the statistics for the countries whose numbers changed -----
# Columns: B=Casos totales, C=Nuevos casos, D=Casos activos, E=Recuperados,
#          F=Casos criticos, G=Muertes hoy, H=Muertes

function Set-RowValues($row, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($row, 2 + $i).Value = $values[$i]
    }
}

# Estados Unidos (row 4)
Set-RowValues 4 @(1410267, 1631, 298641, 1028121, 16473, 80, 83505)

# Mayotte (row 96)
Set-RowValues 96 @(1143, 48, 627, 502, 6, 2, 14)

# Cabo Verde (row 140)
Set-RowValues 140 @(289, 22, 61, 226, 0, 0, 2)

# Santo Tome y Principe (row 145)
Set-RowValues 145 @(220, 12, 4, 210, 0, 1, 6)

# Mauritania (row 211)
Set-RowValues 211 @(9, 0, 6, 1, 0, 1, 2)

# --- 3. Re-sort the data body (rows 4-219) by "Casos totales" descending ---
$dataRange = $ws.Range("A4:H219")
$sortKey = $ws.Range("B4:B219")
$dataRange.Sort($sortKey, 2)
